# This script updates betting-odds values in Sheet1 to match the target
# snapshot described by the commit's XML diff. Each statement sets the
# numeric value of one changed cell; cells not listed are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.62
$ws.Range("G2").Value = 2.98
$ws.Range("H2").Value = 2.8
$ws.Range("I2").Value = 3.45
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 3.5
$ws.Range("P2").Value = 1.65
$ws.Range("Q2").Value = 2.24

# Row 3
$ws.Range("F3").Value = 2.46
$ws.Range("G3").Value = 3.4
$ws.Range("H3").Value = 2.32
$ws.Range("I3").Value = 3.2
$ws.Range("J3").Value = 3.2
$ws.Range("K3").Value = 6.4
$ws.Range("N3").Value = 1.89
$ws.Range("P3").Value = 1.89
$ws.Range("Q3").Value = 1.67
$ws.Range("R3").Value = 1.35
$ws.Range("S3").Value = 2.6

# Row 4
$ws.Range("F4").Value = 13.5
$ws.Range("G4").Value = 19.5
$ws.Range("H4").Value = 1.19
$ws.Range("I4").Value = 1.28
$ws.Range("J4").Value = 7
$ws.Range("K4").Value = 12
$ws.Range("N4").Value = 3.35
$ws.Range("P4").Value = 3.35
$ws.Range("Q4").Value = 1.33
$ws.Range("R4").Value = 1.96
$ws.Range("S4").Value = 1.33
$ws.Range("T4").Value = 1.86
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 4.6
$ws.Range("W4").Value = 1.05

# Row 5
$ws.Range("N5").Value = 5.7
$ws.Range("O5").Value = 1.19
$ws.Range("T5").Value = 1.55
$ws.Range("U5").Value = 2.62
$ws.Range("X5").Value = 25
$ws.Range("Y5").Value = 14.5
$ws.Range("Z5").Value = 16.5
$ws.Range("AA5").Value = 27
$ws.Range("AB5").Value = 20
$ws.Range("AC5").Value = 9.800000000000001
$ws.Range("AD5").Value = 11.5
$ws.Range("AE5").Value = 19.5
$ws.Range("AF5").Value = 29
$ws.Range("AG5").Value = 15.5
$ws.Range("AH5").Value = 15
$ws.Range("AI5").Value = 27
$ws.Range("AJ5").Value = 60
$ws.Range("AK5").Value = 40
$ws.Range("AN5").Value = 23
$ws.Range("AO5").Value = 10.5

# Row 6
$ws.Range("G6").Value = 3.4
$ws.Range("H6").Value = 2.18
$ws.Range("I6").Value = 2.2
$ws.Range("N6").Value = 6.2
$ws.Range("O6").Value = 1.17
$ws.Range("R6").Value = 1.7
$ws.Range("S6").Value = 2.34
$ws.Range("T6").Value = 1.52
$ws.Range("U6").Value = 2.86
$ws.Range("X6").Value = 28
$ws.Range("Y6").Value = 16
$ws.Range("Z6").Value = 18.5
$ws.Range("AA6").Value = 28
$ws.Range("AB6").Value = 21
$ws.Range("AC6").Value = 10.5
$ws.Range("AD6").Value = 12
$ws.Range("AE6").Value = 20
$ws.Range("AF6").Value = 28
$ws.Range("AG6").Value = 15
$ws.Range("AH6").Value = 15
$ws.Range("AI6").Value = 26
$ws.Range("AJ6").Value = 65
$ws.Range("AK6").Value = 32
$ws.Range("AL6").Value = 34
$ws.Range("AM6").Value = 60
$ws.Range("AN6").Value = 19
$ws.Range("AO6").Value = 10

# Row 8
$ws.Range("F8").Value = 3.8
$ws.Range("G8").Value = 4.6
$ws.Range("H8").Value = 1.92
$ws.Range("I8").Value = 2.08
$ws.Range("J8").Value = 3.7
$ws.Range("K8").Value = 4.4
$ws.Range("P8").Value = 2.16
$ws.Range("Q8").Value = 1.69

# Row 9
$ws.Range("F9").Value = 1.8
$ws.Range("G9").Value = 1.99
$ws.Range("H9").Value = 3.85
$ws.Range("I9").Value = 6.4
$ws.Range("J9").Value = 3.45
$ws.Range("K9").Value = 5
$ws.Range("P9").Value = 2.06
$ws.Range("Q9").Value = 1.74

# Row 10
$ws.Range("F10").Value = 2.36
$ws.Range("G10").Value = 2.84
$ws.Range("H10").Value = 2.82
$ws.Range("I10").Value = 3.3
$ws.Range("J10").Value = 3.4
$ws.Range("K10").Value = 4.1
$ws.Range("P10").Value = 2.1
$ws.Range("Q10").Value = 1.72

# Row 11
$ws.Range("F11").Value = 1.42
$ws.Range("G11").Value = 1.53
$ws.Range("H11").Value = 7
$ws.Range("I11").Value = 9.6
$ws.Range("J11").Value = 4.7
$ws.Range("K11").Value = 5.7
$ws.Range("P11").Value = 2.4
$ws.Range("Q11").Value = 1.56

# Row 12
$ws.Range("H12").Value = 4.8
$ws.Range("J12").Value = 3.85
$ws.Range("O12").Value = 1.33
$ws.Range("P12").Value = 1.93
$ws.Range("Q12").Value = 2.02
$ws.Range("X12").Value = 14.5
$ws.Range("Y12").Value = 16.5
$ws.Range("Z12").Value = 38
$ws.Range("AA12").Value = 140
$ws.Range("AB12").Value = 8.800000000000001
$ws.Range("AC12").Value = 8.6
$ws.Range("AD12").Value = 21
$ws.Range("AE12").Value = 75
$ws.Range("AF12").Value = 11
$ws.Range("AG12").Value = 10.5
$ws.Range("AH12").Value = 21
$ws.Range("AI12").Value = 75
$ws.Range("AJ12").Value = 22
$ws.Range("AK12").Value = 20
$ws.Range("AL12").Value = 40
$ws.Range("AM12").Value = 140
$ws.Range("AN12").Value = 12.5
$ws.Range("AO12").Value = 85

# Row 13
$ws.Range("F13").Value = 9.199999999999999
$ws.Range("G13").Value = 9.4
$ws.Range("J13").Value = 6.4
$ws.Range("K13").Value = 6.8
$ws.Range("N13").Value = 8.199999999999999
$ws.Range("O13").Value = 1.12
$ws.Range("P13").Value = 3.4
$ws.Range("R13").Value = 1.98
$ws.Range("T13").Value = 1.67
$ws.Range("X13").Value = 48
$ws.Range("Y13").Value = 15.5
$ws.Range("Z13").Value = 11.5
$ws.Range("AA13").Value = 13.5
$ws.Range("AB13").Value = 50
$ws.Range("AC13").Value = 17
$ws.Range("AD13").Value = 11
$ws.Range("AE13").Value = 13
$ws.Range("AF13").Value = 120
$ws.Range("AG13").Value = 38
$ws.Range("AH13").Value = 24
$ws.Range("AI13").Value = 26
$ws.Range("AK13").Value = 130
$ws.Range("AL13").Value = 85
$ws.Range("AM13").Value = 95
$ws.Range("AO13").Value = 3.7

# Row 14
$ws.Range("F14").Value = 3.8
$ws.Range("G14").Value = 3.85
$ws.Range("H14").Value = 2.02
$ws.Range("I14").Value = 2.04
$ws.Range("J14").Value = 3.95
$ws.Range("O14").Value = 1.23
$ws.Range("R14").Value = 1.54
$ws.Range("S14").Value = 2.68
$ws.Range("T14").Value = 1.62
$ws.Range("U14").Value = 2.52
$ws.Range("X14").Value = 22
$ws.Range("Y14").Value = 12.5
$ws.Range("Z14").Value = 15
$ws.Range("AA14").Value = 25
$ws.Range("AB14").Value = 19
$ws.Range("AE14").Value = 19.5
$ws.Range("AF14").Value = 32
$ws.Range("AG14").Value = 16.5
$ws.Range("AH14").Value = 16
$ws.Range("AI14").Value = 29
$ws.Range("AK14").Value = 40
$ws.Range("AL14").Value = 44
$ws.Range("AM14").Value = 65
$ws.Range("AN14").Value = 32
$ws.Range("AO14").Value = 11

# Row 15
$ws.Range("F15").Value = 1.87
$ws.Range("G15").Value = 1.89
$ws.Range("H15").Value = 4.8
$ws.Range("I15").Value = 5
$ws.Range("J15").Value = 3.75
$ws.Range("M15").Value = 1.08
$ws.Range("N15").Value = 3.55
$ws.Range("P15").Value = 1.89
$ws.Range("Q15").Value = 2.06
$ws.Range("T15").Value = 1.9
$ws.Range("U15").Value = 2
$ws.Range("X15").Value = 14
$ws.Range("Y15").Value = 16.5
$ws.Range("Z15").Value = 36
$ws.Range("AA15").Value = 130
$ws.Range("AB15").Value = 8.6
$ws.Range("AC15").Value = 8.4
$ws.Range("AD15").Value = 19.5
$ws.Range("AE15").Value = 70
$ws.Range("AG15").Value = 10.5
$ws.Range("AH15").Value = 21
$ws.Range("AI15").Value = 75
$ws.Range("AJ15").Value = 22
$ws.Range("AK15").Value = 22
$ws.Range("AL15").Value = 40
$ws.Range("AM15").Value = 140
$ws.Range("AN15").Value = 15
$ws.Range("AO15").Value = 80

# Row 16
$ws.Range("F16").Value = 1.47
$ws.Range("G16").Value = 1.49
$ws.Range("I16").Value = 7.2
$ws.Range("J16").Value = 5.6
$ws.Range("K16").Value = 5.8
$ws.Range("Q16").Value = 1.42
$ws.Range("S16").Value = 2.06
$ws.Range("T16").Value = 1.6
$ws.Range("X16").Value = 40
$ws.Range("Y16").Value = 40
$ws.Range("Z16").Value = 70
$ws.Range("AB16").Value = 15.5
$ws.Range("AC16").Value = 14.5
$ws.Range("AD16").Value = 28
$ws.Range("AE16").Value = 75
$ws.Range("AF16").Value = 13
$ws.Range("AH16").Value = 21
$ws.Range("AI16").Value = 65
$ws.Range("AJ16").Value = 15.5
$ws.Range("AK16").Value = 14
$ws.Range("AL16").Value = 26
$ws.Range("AM16").Value = 70
$ws.Range("AN16").Value = 4.5
$ws.Range("AO16").Value = 55

# Row 17
$ws.Range("H17").Value = 22
$ws.Range("I17").Value = 23
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = 9.199999999999999
$ws.Range("N17").Value = 8
$ws.Range("O17").Value = 1.12
$ws.Range("P17").Value = 3.4
$ws.Range("Q17").Value = 1.38
$ws.Range("R17").Value = 1.98
$ws.Range("S17").Value = 1.93
$ws.Range("T17").Value = 2.16
$ws.Range("X17").Value = 44
$ws.Range("Y17").Value = 990
$ws.Range("AB17").Value = 14
$ws.Range("AC17").Value = 22
$ws.Range("AD17").Value = 990
$ws.Range("AF17").Value = 9.199999999999999
$ws.Range("AG17").Value = 14.5
$ws.Range("AH17").Value = 44
$ws.Range("AJ17").Value = 8.800000000000001
$ws.Range("AK17").Value = 14
$ws.Range("AL17").Value = 46
$ws.Range("AN17").Value = 2.84

# Row 18
$ws.Range("N18").Value = 11
$ws.Range("P18").Value = 4.4
$ws.Range("Q18").Value = 1.26
$ws.Range("S18").Value = 1.66
$ws.Range("T18").Value = 2.06
$ws.Range("AB18").Value = 18.5
$ws.Range("AC18").Value = 50
$ws.Range("AF18").Value = 11.5
$ws.Range("AG18").Value = 16
$ws.Range("AJ18").Value = 9.6
$ws.Range("AK18").Value = 14
$ws.Range("AL18").Value = 46
